# daily auto push: 2025-10-07 13:37 UTC
# Append a new data row (row 75) to the bottom of the log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 75

# Column A ("日付" / date) must stay a plain text string like the other
# date cells in the sheet ("2025/10/07"), not get auto-converted into an
# Excel date serial number. Force text formatting before assigning the
# value, then clear the formatting override so the cell keeps the sheet's
# default (unstyled) look, matching the rest of the data rows.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/07"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = "火"
$ws.Cells.Item($row, 3).Value = 22
$ws.Cells.Item($row, 4).Value = 103
